$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply cell updates per the diff
$ws.Range("D2").Value = "42.292.24"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.248.29"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "247.44"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("D7").Value = "74.81"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").Value = "42.43"
$ws.Range("E10").Value = "  +5.42%  "
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").Value = "7.20"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "14.56"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.858"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.247.45"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "42.119.73"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0986"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "72.01"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "2.25"
$ws.Range("E21").Value = "  +4.30%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "232.18"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "8.89"
$ws.Range("E23").Value = "  +37.84%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "11.25"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  -5.46%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.31"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "169.52"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.10"
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "20.69"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0831"
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "30.26"
$ws.Range("E33").Value = "  -4.89%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.25"
$ws.Range("E35").Value = "  +10.82%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0305"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "13.53"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "2.19"
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "62.19"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").Value = "109.64"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.202"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "8.65"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.101"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "1.16"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "4.17"
$ws.Range("E50").Value = "  -11.67%  "
$ws.Range("E51").Value = "  -1.11%  "
